$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching style of existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data columns I and J, rows 2-16
$values = @{
    2  = @(9, 9)
    3  = @(8, 8)
    4  = @(7, 7)
    5  = @(5, 7)
    6  = @(1, 4)
    7  = @(1, 4)
    8  = @(4, 7)
    9  = @(4, 6)
    10 = @(6, 7)
    11 = @(9, 9)
    12 = @(4, 4)
    13 = @(6, 6)
    14 = @(3, 3)
    15 = @(8, 8)
    16 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
